$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 260, pushing the existing rows 260-301 down to 261-302
$ws.Rows("260:260").Insert()

# Populate the newly inserted row 260 with the new record
$ws.Range("A260").Value = 4
$ws.Range("B260").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C260").Value = "Los Lagos"
$ws.Range("D260").Value = 44951
$ws.Range("E260").Value = 10
$ws.Range("F260").Value = 100112039
$ws.Range("G260").Value = "Ciboulette"
$ws.Range("H260").Value = "Sin especificar"
$ws.Range("I260").Value = "Primera"
$ws.Range("J260").Value = 40
$ws.Range("K260").Value = 6000
$ws.Range("L260").Value = 6000
$ws.Range("M260").Value = 6000
$ws.Range("N260").Value = "$/docena de atados"
$ws.Range("O260").Value = "Provincia de Cautín"
$ws.Range("P260").Value = 2000
$ws.Range("Q260").Value = 3
$ws.Range("R260").Value = "Hortaliza"
